# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets to reflect the latest
# handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 10:22:23"
$wsZhCn.Range("H2").Value = "2016-03-24 10:22:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 10:22:28"
$wsDeDe.Range("H2").Value = "2016-03-24 10:22:55"
